$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '41.813.11'
$ws.Range('E2').Value = '  -1.21%  '
$ws.Range('D3').Value = '2.209.60'
$ws.Range('E3').Value = '  -1.53%  '
$ws.Range('E4').Value = '  +0.36%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '240.39'
$ws.Range('E5').Value = '  -2.28%  '
$ws.Range('E6').Value = '  -1.23%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '72.19'
$ws.Range('E7').Value = '  -4.45%  '
$ws.Range('E8').Value = '  +0.15%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.595'
$ws.Range('E9').Value = '  -3.84%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '41.50'
$ws.Range('E10').Value = '  -5.49%  '
$ws.Range('E11').Value = '  +0.08%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.104'
$ws.Range('E12').Value = '  +1.06%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '6.94'
$ws.Range('E13').Value = '  -4.20%  '
$ws.Range('D14').Value = '2.541.12'
$ws.Range('E14').Value = '  -1.52%  '
$ws.Range('E15').Value = '  -2.98%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.830'
$ws.Range('E16').Value = '  -3.05%  '
$ws.Range('D17').Value = '2.181.78'
$ws.Range('E17').Value = '  -3.02%  '
$ws.Range('D18').Value = '41.700.25'
$ws.Range('E18').Value = '  -1.21%  '
$ws.Range('E19').Value = '  +1.97%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '72.20'
$ws.Range('E20').Value = '  +0.11%  '
$ws.Range('E21').Value = '  -0.84%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '11.08'
$ws.Range('E22').Value = '  +20.07%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '228.33'
$ws.Range('E23').Value = '  -1.42%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.03'
$ws.Range('E24').Value = '  -9.06%  '
$ws.Range('E25').Value = '  +0.04%  '
$ws.Range('E26').Value = '  -1.50%  '
$ws.Range('E27').Value = '  +0.06%  '
$ws.Range('E28').Value = '  -2.72%  '
$ws.Range('E29').Value = '  -0.73%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '166.64'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '20.30'
$ws.Range('E31').Value = '  -1.76%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.55'
$ws.Range('E32').Value = '  +3.87%  '
$ws.Range('E33').Value = '  -4.13%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '29.83'
$ws.Range('E34').Value = '  -2.93%  '
$ws.Range('E35').Value = '  -1.29%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.105'
$ws.Range('E36').Value = '  -12.19%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '4.19'
$ws.Range('E37').Value = '  -7.52%  '
$ws.Range('E38').Value = '  -5.40%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '13.31'
$ws.Range('E39').Value = '  -5.27%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.11'
$ws.Range('E40').Value = '  -3.17%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.59'
$ws.Range('E41').Value = '  -3.46%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '63.67'
$ws.Range('E43').Value = '  -3.45%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '8.61'
$ws.Range('E44').Value = '  -2.15%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '102.62'
$ws.Range('E45').Value = '  -4.67%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0997'
$ws.Range('E46').Value = '  -2.80%  '
$ws.Range('E47').Value = '  -1.37%  '
$ws.Range('E48').Value = '  -2.35%  '
$ws.Range('E49').Value = '  -2.71%  '
$ws.Range('D51').Value = '2.416.74'
$ws.Range('E51').Value = '  -1.62%  '
